$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Install R(3.5.1) into directory: C:/Program Files/R/R-3.5.1"
#    -> append " (or another folder)" at the end of the paragraph.
# ------------------------------------------------------------------
$find1 = $d.Content
$found1 = $find1.Find.Execute("Install R(3.5.1) into directory: C:/Program Files/R/R-3.5.1")
if ($found1) {
    $insertPoint = $d.Range($find1.End, $find1.End)
    $insertPoint.InsertAfter(" (or another folder)")
}

# ------------------------------------------------------------------
# 2) "C:/Program Files/R/R-3.2.5/bin/Rscript" ->
#    "C:/Program Files/R/R-3.5.1/bin/Rscript (or the path you specified)"
#    with the _GoBack bookmark relocated right after that new text
#    (just before the closing curly quote).
# ------------------------------------------------------------------
$find2 = $d.Content
$found2 = $find2.Find.Execute("C:/Program Files/R/R-3.2.5/bin/Rscript")
if ($found2) {
    $start = $find2.Start
    $find2.Text = "C:/Program Files/R/R-3.5.1/bin/Rscript"
    $afterPath = $start + ("C:/Program Files/R/R-3.5.1/bin/Rscript").Length

    $insertPoint2 = $d.Range($afterPath, $afterPath)
    $insertPoint2.InsertAfter(" (or the path you specified)")
    $afterNote = $afterPath + (" (or the path you specified)").Length

    # Move the _GoBack bookmark here (Bookmarks.Add with an existing
    # name re-defines/moves it, removing the old occurrence).
    $bmRange = $d.Range($afterNote, $afterNote)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

# ------------------------------------------------------------------
# 3) Remove the stray <w:lastRenderedPageBreak/> that sits in front of
#    the Google Drive link run (re-typing the run's text drops the
#    stale rendering bookkeeping element).
# ------------------------------------------------------------------
$find3 = $d.Content
$found3 = $find3.Find.Execute("https://drive.google.com/drive/folders/1bD77wI0-nT-j5qUkLGuM6S5DgmtG1Jb0?usp=sharing")
if ($found3) {
    $linkStart = $find3.Start
    $linkText = $find3.Text
    $find3.Text = ""
    $reinsert = $d.Range($linkStart, $linkStart)
    $reinsert.InsertAfter($linkText)
}
